# ===================================================================
# Applies the recalculation update described by the commit:
#   mdl.medium_heat[t] == steam_heat_load[t] + mdl.DRHeatLoad[t] if t in
#   peak else heat cannot be wasted during off-peak shaving periods
#
# The model horizon grew from 92 to 96 periods (t = 0..95), so each of
# the six result sheets gains four new data rows (94-97) and a handful
# of pre-existing cells receive updated (re-solved) numeric values.
# ===================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# sheet1 (Worksheets.Item(1))
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item(1)

# Updated values coming from the re-solved model
$ws.Range("G2").Value = -0.888888888888914
$ws.Range("G3").Value = -0.888888888888942
$ws.Range("G10").Value = -50.7421147707863
$ws.Range("C14").Value = 7309.46492238139
$ws.Range("C36").Value = 3938.43558594186
$ws.Range("E36").Value = 0
$ws.Range("C37").Value = 4216.54426086744
$ws.Range("E37").Value = 0
$ws.Range("C38").Value = 4271.42598422791
$ws.Range("E38").Value = 0
$ws.Range("C39").Value = 4245.62911115581
$ws.Range("E39").Value = 0
$ws.Range("E70").Value = 625
$ws.Range("E71").Value = 0
$ws.Range("C72").Value = 4460.85626411512
$ws.Range("H72").Value = -0
$ws.Range("C73").Value = 4858.36459960275
$ws.Range("F73").Value = 1000
$ws.Range("H73").Value = -290.342075149028
$ws.Range("C74").Value = 5607.77363355135
$ws.Range("H74").Value = -915.342075149028
$ws.Range("C75").Value = 5329.08124735465
$ws.Range("H75").Value = -625
$ws.Range("C81").Value = 3938.03386091454
$ws.Range("C82").Value = 4042.58150855362
$ws.Range("I82").Value = -264.03070181297

# Append new rows 94-97: copy formatting of column A from
# row 93 (bold/centered/bordered index style) so the new index cells
# match the existing table formatting, then fill in the values.
$ws.Range("A93").Copy() | Out-Null
$ws.Range("A94:A97").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A94").Value = 92
$ws.Range("B94").Value = -5669.4340791
$ws.Range("C94").Value = 5058.57094078911
$ws.Range("D94").Value = 0
$ws.Range("E94").Value = 0
$ws.Range("F94").Value = 1000
$ws.Range("G94").Value = -0
$ws.Range("H94").Value = -13.6842105263158
$ws.Range("I94").Value = -357.57395348838
$ws.Range("J94").Value = 0.65550003
$ws.Range("A95").Value = 93
$ws.Range("B95").Value = -5769.46635825
$ws.Range("C95").Value = 5153.16591080423
$ws.Range("D95").Value = 0
$ws.Range("E95").Value = 0
$ws.Range("F95").Value = 1000
$ws.Range("G95").Value = -0
$ws.Range("H95").Value = -13.6842105263158
$ws.Range("I95").Value = -351.945614883727
$ws.Range("J95").Value = 0.65550003
$ws.Range("A96").Value = 94
$ws.Range("B96").Value = -5819.4986374
$ws.Range("C96").Value = 5196.36673131237
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("F96").Value = 1000
$ws.Range("G96").Value = -0
$ws.Range("H96").Value = -13.6842105263158
$ws.Range("I96").Value = -344.838072558146
$ws.Range("J96").Value = 0.65550003
$ws.Range("A97").Value = 95
$ws.Range("B97").Value = -5869.4986374
$ws.Range("C97").Value = 5243.0422494147
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("F97").Value = 1000
$ws.Range("G97").Value = -0
$ws.Range("H97").Value = -13.6842105263158
$ws.Range("I97").Value = -341.451460465124
$ws.Range("J97").Value = 0.65550003


# ---------------------------------------------------------------
# sheet2 (Worksheets.Item(2))
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item(2)

# Updated values coming from the re-solved model
$ws.Range("C2").Value = [double]"2.220446049250313e-15"
$ws.Range("C3").Value = [double]"5.773159728050814e-15"
$ws.Range("C10").Value = [double]"-1.13686837721616e-13"
$ws.Range("B82").Value = 1135.332017795771
$ws.Range("C82").Value = 416.379982204229

# Append new rows 94-97: copy formatting of column A from
# row 93 (bold/centered/bordered index style) so the new index cells
# match the existing table formatting, then fill in the values.
$ws.Range("A93").Copy() | Out-Null
$ws.Range("A94:A97").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A94").Value = 92
$ws.Range("B94").Value = 1537.568000000034
$ws.Range("C94").Value = 0
$ws.Range("D94").Value = 0
$ws.Range("E94").Value = -1537.568000000032
$ws.Range("F94").Value = 0.65550003
$ws.Range("A95").Value = 93
$ws.Range("B95").Value = 1513.366144000026
$ws.Range("C95").Value = 0
$ws.Range("D95").Value = 18.89785600000944
$ws.Range("E95").Value = -1532.264000000036
$ws.Range("F95").Value = 0.65550003
$ws.Range("A96").Value = 94
$ws.Range("B96").Value = 1482.803712000028
$ws.Range("C96").Value = 0
$ws.Range("D96").Value = 44.15628800001105
$ws.Range("E96").Value = -1526.96000000004
$ws.Range("F96").Value = 0.65550003
$ws.Range("A97").Value = 95
$ws.Range("B97").Value = 1468.241280000033
$ws.Range("C97").Value = 0
$ws.Range("D97").Value = 53.41472000001264
$ws.Range("E97").Value = -1521.656000000044
$ws.Range("F97").Value = 0.65550003


# ---------------------------------------------------------------
# sheet3 (Worksheets.Item(3))
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item(3)

# Updated values coming from the re-solved model
$ws.Range("B2").Value = [double]"2.220446049250313e-15"
$ws.Range("B3").Value = [double]"5.773159728050814e-15"
$ws.Range("B10").Value = [double]"-1.13686837721616e-13"
$ws.Range("B82").Value = 416.379982204229
$ws.Range("C82").Value = 1135.332017795771

# Append new rows 94-97: copy formatting of column A from
# row 93 (bold/centered/bordered index style) so the new index cells
# match the existing table formatting, then fill in the values.
$ws.Range("A93").Copy() | Out-Null
$ws.Range("A94:A97").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A94").Value = 92
$ws.Range("B94").Value = 0
$ws.Range("C94").Value = 1537.568000000034
$ws.Range("D94").Value = 0
$ws.Range("E94").Value = 1537.568000000032
$ws.Range("F94").Value = 0.65550003
$ws.Range("A95").Value = 93
$ws.Range("B95").Value = 0
$ws.Range("C95").Value = 1513.366144000026
$ws.Range("D95").Value = 18.89785600000944
$ws.Range("E95").Value = 1532.264000000036
$ws.Range("F95").Value = 0.65550003
$ws.Range("A96").Value = 94
$ws.Range("B96").Value = 0
$ws.Range("C96").Value = 1482.803712000028
$ws.Range("D96").Value = 44.15628800001105
$ws.Range("E96").Value = 1526.96000000004
$ws.Range("F96").Value = 0.65550003
$ws.Range("A97").Value = 95
$ws.Range("B97").Value = 0
$ws.Range("C97").Value = 1468.241280000033
$ws.Range("D97").Value = 53.41472000001264
$ws.Range("E97").Value = 1521.656000000044
$ws.Range("F97").Value = 0.65550003


# ---------------------------------------------------------------
# sheet4 (Worksheets.Item(4))
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item(4)

# Updated values coming from the re-solved model
$ws.Range("B73").Value = 4133.02378181818
$ws.Range("C73").Value = 2030.30303030303

# Append new rows 94-97: copy formatting of column A from
# row 93 (bold/centered/bordered index style) so the new index cells
# match the existing table formatting, then fill in the values.
$ws.Range("A93").Copy() | Out-Null
$ws.Range("A94:A97").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A94").Value = 92
$ws.Range("B94").Value = 2108.57530181822
$ws.Range("C94").Value = 2030.30303030303
$ws.Range("D94").Value = 0
$ws.Range("E94").Value = -3318.80640000005
$ws.Range("F94").Value = -3318.80640000005
$ws.Range("G94").Value = 0.65550003
$ws.Range("A95").Value = 93
$ws.Range("B95").Value = 2233.32978181824
$ws.Range("C95").Value = 2030.30303030303
$ws.Range("D95").Value = 23.6223200000118
$ws.Range("E95").Value = -3451.51160000006
$ws.Range("F95").Value = -3451.51160000006
$ws.Range("G95").Value = 0.65550003
$ws.Range("A96").Value = 94
$ws.Range("B96").Value = 2366.03498181825
$ws.Range("C96").Value = 2030.30303030303
$ws.Range("D96").Value = 55.19536000001381
$ws.Range("E96").Value = -3584.21680000007
$ws.Range("F96").Value = -3584.21680000007
$ws.Range("G96").Value = 0.65550003
$ws.Range("A97").Value = 95
$ws.Range("B97").Value = 2398.74018181826
$ws.Range("C97").Value = 2030.30303030303
$ws.Range("D97").Value = 66.7684000000158
$ws.Range("E97").Value = -3616.92200000008
$ws.Range("F97").Value = -3616.92200000008
$ws.Range("G97").Value = 0.65550003


# ---------------------------------------------------------------
# sheet5 (Worksheets.Item(5))
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item(5)

# Updated values coming from the re-solved model
$ws.Range("B73").Value = 4133.02378181818
$ws.Range("C73").Value = 1218.181818181818
$ws.Range("D73").Value = 812.121212121212
$ws.Range("F73").Value = 4133.02378181818

# Append new rows 94-97: copy formatting of column A from
# row 93 (bold/centered/bordered index style) so the new index cells
# match the existing table formatting, then fill in the values.
$ws.Range("A93").Copy() | Out-Null
$ws.Range("A94:A97").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A94").Value = 92
$ws.Range("B94").Value = 2108.57530181822
$ws.Range("C94").Value = 1218.181818181818
$ws.Range("D94").Value = 812.121212121212
$ws.Range("E94").Value = 663.7612800000101
$ws.Range("F94").Value = 2108.57530181822
$ws.Range("G94").Value = -3318.80640000005
$ws.Range("H94").Value = -3318.80640000005
$ws.Range("I94").Value = -0
$ws.Range("J94").Value = 0
$ws.Range("A95").Value = 93
$ws.Range("B95").Value = 2233.32978181824
$ws.Range("C95").Value = 1218.181818181818
$ws.Range("D95").Value = 812.121212121212
$ws.Range("E95").Value = 690.3023200000121
$ws.Range("F95").Value = 2233.32978181824
$ws.Range("G95").Value = -3451.51160000006
$ws.Range("H95").Value = -3451.51160000006
$ws.Range("I95").Value = -23.6223200000118
$ws.Range("J95").Value = 18.89785600000944
$ws.Range("A96").Value = 94
$ws.Range("B96").Value = 2366.03498181825
$ws.Range("C96").Value = 1218.181818181818
$ws.Range("D96").Value = 812.121212121212
$ws.Range("E96").Value = 716.8433600000141
$ws.Range("F96").Value = 2366.03498181825
$ws.Range("G96").Value = -3584.21680000007
$ws.Range("H96").Value = -3584.21680000007
$ws.Range("I96").Value = -55.19536000001381
$ws.Range("J96").Value = 44.15628800001105
$ws.Range("A97").Value = 95
$ws.Range("B97").Value = 2398.74018181826
$ws.Range("C97").Value = 1218.181818181818
$ws.Range("D97").Value = 812.121212121212
$ws.Range("E97").Value = 723.3844000000161
$ws.Range("F97").Value = 2398.74018181826
$ws.Range("G97").Value = -3616.92200000008
$ws.Range("H97").Value = -3616.92200000008
$ws.Range("I97").Value = -66.7684000000158
$ws.Range("J97").Value = 53.41472000001264


# ---------------------------------------------------------------
# sheet6 (Worksheets.Item(6))
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item(6)

# Updated values coming from the re-solved model
$ws.Range("C73").Value = 812.121212121212

# Append new rows 94-97: copy formatting of column A from
# row 93 (bold/centered/bordered index style) so the new index cells
# match the existing table formatting, then fill in the values.
$ws.Range("A93").Copy() | Out-Null
$ws.Range("A94:A97").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A94").Value = 92
$ws.Range("B94").Value = 663.7612800000101
$ws.Range("C94").Value = 812.121212121212
$ws.Range("D94").Value = -0
$ws.Range("E94").Value = 0.65550003
$ws.Range("A95").Value = 93
$ws.Range("B95").Value = 690.3023200000121
$ws.Range("C95").Value = 812.121212121212
$ws.Range("D95").Value = -23.6223200000118
$ws.Range("E95").Value = 0.65550003
$ws.Range("A96").Value = 94
$ws.Range("B96").Value = 716.8433600000141
$ws.Range("C96").Value = 812.121212121212
$ws.Range("D96").Value = -55.19536000001381
$ws.Range("E96").Value = 0.65550003
$ws.Range("A97").Value = 95
$ws.Range("B97").Value = 723.3844000000161
$ws.Range("C97").Value = 812.121212121212
$ws.Range("D97").Value = -66.7684000000158
$ws.Range("E97").Value = 0.65550003

